$wb = $excel.ActiveWorkbook

# Update "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# on the zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-11 18:46:43"
$wsZh.Range("H2").Value = "2016-03-11 18:47:14"

# Update the same columns on the de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-11 18:46:47"
$wsDe.Range("H2").Value = "2016-03-11 18:47:20"
